# Practiced Cambridge 7 Test1 Listening.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the Writing (I) scores for the already-logged tests (rows 20-24) ---
# These were placeholder 1.1 values; replace them with the real scores, which
# also recalculates the shared Overall (K) formula for each row.
$ws.Range("I20").Value2 = 5.5
$ws.Range("I21").Value2 = 5.5
$ws.Range("I22").Value2 = 6.5
$ws.Range("I23").Value2 = 6
$ws.Range("I24").Value2 = 6.5

# --- Add a new practice-test entry in row 25 (Cambridge 7 Test 1 Listening) ---
# First clone the formatting from the row above (column-by-column, skipping G
# which carries no explicit style) so the new row's cell styles/borders match.
$ws.Range("C24:F24").Copy()
$ws.Range("C25:F25").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H24:K24").Copy()
$ws.Range("H25:K25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Now fill in the new row's data.
$ws.Range("C25").Value2 = 45488   # 15-Jul-2024
$ws.Range("C25").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("D25").Value = "IELTS7_Test1"
$ws.Range("E25").Value2 = 36
$ws.Range("F25").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20,MATCH([@Lis_Mark],Sheet2!$D$5:$D$20,1)),"No Grade")'

# --- Restore the author's final view state (scroll position + selection) ---
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I25").Select()

$wb.Save()
